$wb = $excel.ActiveWorkbook
$ws5 = $wb.Worksheets.Item("unified_perf_summary")
$ws8 = $wb.Worksheets.Item("Normalization")

# --- unified_perf_summary sheet ---
$ws5.Range("P2").Value = 1.75
$ws5.Range("Q2").Value = 0.6428571428571429
$ws5.Range("S2").Value = 0.05628963269412802
$ws5.Range("T2").Value = 0.001008475034056347
$ws5.Range("AC2").Value = 0.05683579419860258
$ws5.Range("AD2").Value = 0.05512587657870418
$ws5.Range("AE2").Value = 0.0569072273050773

$ws5.Range("P3").Value = 1.5625
$ws5.Range("Q3").Value = 0.72
$ws5.Range("S3").Value = 0.06680599771272798
$ws5.Range("T3").Value = 0.00219379562352543
$ws5.Range("AC3").Value = 0.06661193893554086
$ws5.Range("AD3").Value = 0.06471567821944493
$ws5.Range("AE3").Value = 0.06909037598319813

$ws5.Range("P5").Value = 1.50006103515625
$ws5.Range("Q5").Value = 0.749969483663588
$ws5.Range("S5").Value = 0.08158110778877121
$ws5.Range("T5").Value = 0.008789629375982646
$ws5.Range("AC5").Value = 0.07823573877352762
$ws5.Range("AD5").Value = 0.07495536086744072
$ws5.Range("AE5").Value = 0.0915522237253453

$ws5.Range("P12").Value = 1.50006103515625
$ws5.Range("Q12").Value = 0.749969483663588
$ws5.Range("S12").Value = 0.1160727043946573
$ws5.Range("T12").Value = 0.02471211817896858
$ws5.Range("AC12").Value = 0.1204695790575916
$ws5.Range("AD12").Value = 0.08945727697861705
$ws5.Range("AE12").Value = 0.1382912571477634

$ws5.Range("P14").Value = 1.500244140625
$ws5.Range("Q14").Value = 0.7498779495524817
$ws5.Range("S14").Value = 0.13861462663244
$ws5.Range("T14").Value = 0.001021216688976221
$ws5.Range("AC14").Value = 0.138284391793287
$ws5.Range("AD14").Value = 0.1377993909324209
$ws5.Range("AE14").Value = 0.139760097171612

$ws5.Range("P20").Value = 1.501953125
$ws5.Range("Q20").Value = 0.7490247074122237
$ws5.Range("S20").Value = 0.241821376849506
$ws5.Range("T20").Value = 0.0297823893140027
$ws5.Range("AC20").Value = 0.2488365820089493
$ws5.Range("AD20").Value = 0.2091576276506063
$ws5.Range("AE20").Value = 0.2674699208889626

$ws5.Range("P23").Value = 1.00006103515625
$ws5.Range("Q23").Value = 0.5000152578577968
$ws5.Range("S23").Value = 0.2067490101745708
$ws5.Range("T23").Value = 0.01852446773236728
$ws5.Range("AC23").Value = 0.2094416539886874
$ws5.Range("AD23").Value = 0.1870255786815292
$ws5.Range("AE23").Value = 0.2237797978534959

$ws5.Range("P25").Value = 1.00048828125
$ws5.Range("Q25").Value = 0.5001220107369448
$ws5.Range("S25").Value = 0.2371125493316399
$ws5.Range("T25").Value = 0.02942084384243977
$ws5.Range("AC25").Value = 0.2447632973342447
$ws5.Range("AD25").Value = 0.2046221165714285
$ws5.Range("AE25").Value = 0.2619522340892465

# --- Normalization sheet ---
$ws8.Range("L3").Value = 1.75
$ws8.Range("M3").Value = 0.6428571428571429
$ws8.Range("N3").Value = 0.05628963269412802
$ws8.Range("O3").Value = 0.05683579419860258
$ws8.Range("P3").Value = 0.001008475034056347
$ws8.Range("Q3").Value = 0.05512587657870418
$ws8.Range("R3").Value = 0.0569072273050773

$ws8.Range("L9").Value = 1.5625
$ws8.Range("M9").Value = 0.72
$ws8.Range("N9").Value = 0.06680599771272798
$ws8.Range("O9").Value = 0.06661193893554086
$ws8.Range("P9").Value = 0.00219379562352543
$ws8.Range("Q9").Value = 0.06471567821944493
$ws8.Range("R9").Value = 0.06909037598319813

$ws8.Range("L12").Value = 1.50006103515625
$ws8.Range("M12").Value = 0.749969483663588
$ws8.Range("N12").Value = 0.08158110778877121
$ws8.Range("O12").Value = 0.07823573877352762
$ws8.Range("P12").Value = 0.008789629375982646
$ws8.Range("Q12").Value = 0.07495536086744072
$ws8.Range("R12").Value = 0.0915522237253453

$ws8.Range("L16").Value = 1.50006103515625
$ws8.Range("M16").Value = 0.749969483663588
$ws8.Range("N16").Value = 0.1160727043946573
$ws8.Range("O16").Value = 0.1204695790575916
$ws8.Range("P16").Value = 0.02471211817896858
$ws8.Range("Q16").Value = 0.08945727697861705
$ws8.Range("R16").Value = 0.1382912571477634

$ws8.Range("L17").Value = 1.500244140625
$ws8.Range("M17").Value = 0.7498779495524817
$ws8.Range("N17").Value = 0.13861462663244
$ws8.Range("O17").Value = 0.138284391793287
$ws8.Range("P17").Value = 0.001021216688976221
$ws8.Range("Q17").Value = 0.1377993909324209
$ws8.Range("R17").Value = 0.139760097171612

$ws8.Range("L19").Value = 1.501953125
$ws8.Range("M19").Value = 0.7490247074122237
$ws8.Range("N19").Value = 0.241821376849506
$ws8.Range("O19").Value = 0.2488365820089493
$ws8.Range("P19").Value = 0.0297823893140027
$ws8.Range("Q19").Value = 0.2091576276506063
$ws8.Range("R19").Value = 0.2674699208889626

$ws8.Range("L21").Value = 1.00006103515625
$ws8.Range("M21").Value = 0.5000152578577968
$ws8.Range("N21").Value = 0.2067490101745708
$ws8.Range("O21").Value = 0.2094416539886874
$ws8.Range("P21").Value = 0.01852446773236728
$ws8.Range("Q21").Value = 0.1870255786815292
$ws8.Range("R21").Value = 0.2237797978534959

$ws8.Range("L22").Value = 1.00048828125
$ws8.Range("M22").Value = 0.5001220107369448
$ws8.Range("N22").Value = 0.2371125493316399
$ws8.Range("O22").Value = 0.2447632973342447
$ws8.Range("P22").Value = 0.02942084384243977
$ws8.Range("Q22").Value = 0.2046221165714285
$ws8.Range("R22").Value = 0.2619522340892465
